$wb = $excel.ActiveWorkbook

# Row -> new F-column value, applied identically to both "展览" and "全部类型" sheets
$updates = @{
    4  = 43
    5  = 2611
    7  = 133
    9  = 1353
    11 = 57
    12 = 14
    13 = 1174
    14 = 347
    15 = 324
    16 = 34
    17 = 29
    21 = 2446
    22 = 27
    23 = 280
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
